# Insert a new weekly record at row 432 in the "Zanahoria" data sheet.
# This shifts all the existing data rows (old 432..490) down by one row
# (new 433..491), and the new row 432 receives this week's fresh values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row before row 432 - pushes row 432..490 down to 433..491
$ws.Rows.Item(432).Insert()

# Populate the newly inserted row 432 with the new weekly record
$ws.Cells.Item(432, 1).Value2  = 7
$ws.Cells.Item(432, 2).Value2  = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(432, 3).Value2  = "Ñuble"
$ws.Cells.Item(432, 4).Value2  = 45131
$ws.Cells.Item(432, 5).Value2  = 16
$ws.Cells.Item(432, 6).Value2  = 100114013
$ws.Cells.Item(432, 7).Value2  = "Zanahoria"
$ws.Cells.Item(432, 8).Value2  = "Sin especificar"
$ws.Cells.Item(432, 9).Value2  = "Primera"
$ws.Cells.Item(432, 10).Value2 = 120
$ws.Cells.Item(432, 11).Value2 = 6000
$ws.Cells.Item(432, 12).Value2 = 6000
$ws.Cells.Item(432, 13).Value2 = 6000
$ws.Cells.Item(432, 14).Value2 = "$/saco 20 kilos"
$ws.Cells.Item(432, 15).Value2 = "Provincia de Diguillín"
$ws.Cells.Item(432, 16).Value2 = 300
$ws.Cells.Item(432, 17).Value2 = 20
$ws.Cells.Item(432, 18).Value2 = "Hortaliza"

# Make sure the new row's date cell uses the same date number format as the
# rest of the "Fecha" column (style index copied automatically by Insert,
# but set explicitly to be safe).
$ws.Cells.Item(432, 4).NumberFormat = $ws.Cells.Item(433, 4).NumberFormat
